$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 10-17 (rows shifted down to make room for a late "PACTO" trade
# entered as row 10, which was missing from the original report) and append the new
# rows 18-25 that the corrected report includes.

# Row 10
$ws.Cells.Item(10, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(10, 2).Value = 45540
$ws.Cells.Item(10, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(10, 3).Value = 45539
$ws.Cells.Item(10, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(10, 4).Value = 0.48
$ws.Cells.Item(10, 5).Value = 22003520
$ws.Cells.Item(10, 6).Value = 22000000
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = "PACTO"
$ws.Cells.Item(10, 9).Value = "COMPRA"
$ws.Cells.Item(10, 10).Value = "PACTO"

# Row 11
$ws.Cells.Item(11, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(11, 2).Value = 45539
$ws.Cells.Item(11, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(11, 3).Value = 45539
$ws.Cells.Item(11, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(11, 4).Value = 5.4
$ws.Cells.Item(11, 5).Value = 500000000
$ws.Cells.Item(11, 6).Value = 522623709
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = "BTP0581029"
$ws.Cells.Item(11, 9).Value = "VENTA"
$ws.Cells.Item(11, 10).Value = "RENTA FIJA"

# Row 12
$ws.Cells.Item(12, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(12, 2).Value = 45539
$ws.Cells.Item(12, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(12, 3).Value = 45539
$ws.Cells.Item(12, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(12, 4).Value = 5.87
$ws.Cells.Item(12, 5).Value = 2000000000
$ws.Cells.Item(12, 6).Value = 1900159321
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = "BBCIO21220"
$ws.Cells.Item(12, 9).Value = "VENTA"
$ws.Cells.Item(12, 10).Value = "RENTA FIJA"

# Row 13
$ws.Cells.Item(13, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(13, 2).Value = 45539
$ws.Cells.Item(13, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(13, 3).Value = 45539
$ws.Cells.Item(13, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(13, 4).Value = 3.6
$ws.Cells.Item(13, 5).Value = 5000
$ws.Cells.Item(13, 6).Value = 192289952
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = "BCAJBF0322"
$ws.Cells.Item(13, 9).Value = "VENTA"
$ws.Cells.Item(13, 10).Value = "RENTA FIJA"

# Row 14
$ws.Cells.Item(14, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(14, 2).Value = 45539
$ws.Cells.Item(14, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(14, 3).Value = 45539
$ws.Cells.Item(14, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(14, 4).Value = 5.85
$ws.Cells.Item(14, 5).Value = 500000000
$ws.Cells.Item(14, 6).Value = 488642906
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = "BCHIDU0716"
$ws.Cells.Item(14, 9).Value = "VENTA"
$ws.Cells.Item(14, 10).Value = "RENTA FIJA"

# Row 15
$ws.Cells.Item(15, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(15, 2).Value = 45539
$ws.Cells.Item(15, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(15, 3).Value = 45539
$ws.Cells.Item(15, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(15, 4).Value = 2.36
$ws.Cells.Item(15, 5).Value = 23000
$ws.Cells.Item(15, 6).Value = 900457256
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = "BCHIAZ0613"
$ws.Cells.Item(15, 9).Value = "COMPRA"
$ws.Cells.Item(15, 10).Value = "RENTA FIJA"

# Row 16
$ws.Cells.Item(16, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(16, 2).Value = 45539
$ws.Cells.Item(16, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(16, 3).Value = 45539
$ws.Cells.Item(16, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(16, 4).Value = 2.06
$ws.Cells.Item(16, 5).Value = 1000
$ws.Cells.Item(16, 6).Value = 38066280
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = "BSECK70915"
$ws.Cells.Item(16, 9).Value = "COMPRA"
$ws.Cells.Item(16, 10).Value = "RENTA FIJA"

# Row 17
$ws.Cells.Item(17, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(17, 2).Value = 45539
$ws.Cells.Item(17, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(17, 3).Value = 45539
$ws.Cells.Item(17, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(17, 4).Value = 2.05
$ws.Cells.Item(17, 5).Value = 2000
$ws.Cells.Item(17, 6).Value = 76139894
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = "BSECK70915"
$ws.Cells.Item(17, 9).Value = "COMPRA"
$ws.Cells.Item(17, 10).Value = "RENTA FIJA"

# Row 18
$ws.Cells.Item(18, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(18, 2).Value = 45539
$ws.Cells.Item(18, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(18, 3).Value = 45539
$ws.Cells.Item(18, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(18, 4).Value = 2.05
$ws.Cells.Item(18, 5).Value = 1000
$ws.Cells.Item(18, 6).Value = 38059258
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = "BSECK70915"
$ws.Cells.Item(18, 9).Value = "COMPRA"
$ws.Cells.Item(18, 10).Value = "RENTA FIJA"

# Row 19
$ws.Cells.Item(19, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(19, 2).Value = 45539
$ws.Cells.Item(19, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(19, 3).Value = 45539
$ws.Cells.Item(19, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(19, 4).Value = 2.2
$ws.Cells.Item(19, 5).Value = 30000
$ws.Cells.Item(19, 6).Value = 1132814474
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = "BBNSAN0918"
$ws.Cells.Item(19, 9).Value = "COMPRA"
$ws.Cells.Item(19, 10).Value = "RENTA FIJA"

# Row 20
$ws.Cells.Item(20, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(20, 2).Value = 45539
$ws.Cells.Item(20, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(20, 3).Value = 45539
$ws.Cells.Item(20, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(20, 4).Value = 1.55
$ws.Cells.Item(20, 5).Value = 30000
$ws.Cells.Item(20, 6).Value = 1169147129
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = "BCODE-B"
$ws.Cells.Item(20, 9).Value = "COMPRA"
$ws.Cells.Item(20, 10).Value = "RENTA FIJA"

# Row 21
$ws.Cells.Item(21, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(21, 2).Value = 45539
$ws.Cells.Item(21, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(21, 3).Value = 45539
$ws.Cells.Item(21, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 18044822.5
$ws.Cells.Item(21, 6).Value = 18037764
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = "SQM-B"
$ws.Cells.Item(21, 9).Value = "VENTA"
$ws.Cells.Item(21, 10).Value = "SIMULTANEA"

# Row 22
$ws.Cells.Item(22, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(22, 2).Value = 45539
$ws.Cells.Item(22, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(22, 3).Value = 45539
$ws.Cells.Item(22, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(22, 4).Value = 16425
$ws.Cells.Item(22, 5).Value = 919
$ws.Cells.Item(22, 6).Value = 15094575
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = "CFINHRFLA"
$ws.Cells.Item(22, 9).Value = "COMPRA"
$ws.Cells.Item(22, 10).Value = "RENTA VARIABLE"

# Row 23
$ws.Cells.Item(23, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(23, 2).Value = 45539
$ws.Cells.Item(23, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(23, 3).Value = 45539
$ws.Cells.Item(23, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(23, 4).Value = 16425.05
$ws.Cells.Item(23, 5).Value = 147
$ws.Cells.Item(23, 6).Value = 2414483
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = "CFINHRFLA"
$ws.Cells.Item(23, 9).Value = "COMPRA"
$ws.Cells.Item(23, 10).Value = "RENTA VARIABLE"

# Row 24
$ws.Cells.Item(24, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(24, 2).Value = 45539
$ws.Cells.Item(24, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(24, 3).Value = 45539
$ws.Cells.Item(24, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(24, 4).Value = 16425.05
$ws.Cells.Item(24, 5).Value = 1091
$ws.Cells.Item(24, 6).Value = 17919730
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = "CFINHRFLA"
$ws.Cells.Item(24, 9).Value = "VENTA"
$ws.Cells.Item(24, 10).Value = "RENTA VARIABLE"

# Row 25
$ws.Cells.Item(25, 1).Value = "FONDO DE INVERSION NEVASA AHORRO"
$ws.Cells.Item(25, 2).Value = 45539
$ws.Cells.Item(25, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(25, 3).Value = 45539
$ws.Cells.Item(25, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(25, 4).Value = 1183.88
$ws.Cells.Item(25, 5).Value = 2535
$ws.Cells.Item(25, 6).Value = 3001136
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = "CFINVRFLI"
$ws.Cells.Item(25, 9).Value = "COMPRA"
$ws.Cells.Item(25, 10).Value = "RENTA VARIABLE"

